$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(416, 1).Value = 415
$ws.Cells.Item(416, 2).Value = 'Pós OS'
$ws.Cells.Item(416, 3).Value = 8005291277
$ws.Cells.Item(416, 4).Value = 46064.658518518518
$ws.Cells.Item(416, 5).Value = 'FRQ_ECO_SP_OSASCO'
$ws.Cells.Item(416, 6).Value = 'Detratores'
$ws.Cells.Item(416, 7).Value = 'A visita NÃO ACONTECEU!!!! Apesar de ter sido marcada a visita para o período da TARDE, o técnico veio no período da MANHÃ!!!! Eu não podia atende-lo pois estava no banho me preparando para uma aula, e o meu Marly estava se arrumando para a fisioterapia que ele tinha em seguida. Me ligaram do atendimento Brastemp, eu expliquei a situação e a moça que me informou que ele viria às 13hs, Mas ele NÃO apareceu!!!!!'
$ws.Cells.Item(416, 8).Value = 'Campo'
$ws.Cells.Item(416, 9).Value = 'Fora do período agendado'

$ws.Cells.Item(417, 1).Value = 416
$ws.Cells.Item(417, 2).Value = 'Pós OS'
$ws.Cells.Item(417, 3).Value = 8005284760
$ws.Cells.Item(417, 4).Value = 46065.372743055559
$ws.Cells.Item(417, 5).Value = 'FRQ_ECO_PE_RECIFE'
$ws.Cells.Item(417, 6).Value = 'Detratores'
$ws.Cells.Item(417, 7).Value = 'Até hoje não foi resolvido, continuamos tomando água quente'
$ws.Cells.Item(417, 8).Value = 'Qualidade do Produto'
$ws.Cells.Item(417, 9).Value = 'Não gela'

$ws.Cells.Item(418, 1).Value = 417
$ws.Cells.Item(418, 2).Value = 'Pós OS'
$ws.Cells.Item(418, 3).Value = 8005267408
$ws.Cells.Item(418, 4).Value = 46065.383090277777
$ws.Cells.Item(418, 5).Value = 'FRQ_ECO_SP_OSASCO'
$ws.Cells.Item(418, 6).Value = 'Detratores'
$ws.Cells.Item(418, 7).Value = 'dei zero para o atendimento de solicitar o técnico, pois demorou quase 2 meses, para poder vir fazer a manutenção preventiva. quanto o tecnico, fez o seu trabalho em 20 mn e foi embora.'
$ws.Cells.Item(418, 8).Value = 'Capacidade'
$ws.Cells.Item(418, 9).Value = 'Data Distante'

$ws.Cells.Item(419, 1).Value = 418
$ws.Cells.Item(419, 2).Value = 'Pós OS'
$ws.Cells.Item(419, 3).Value = 8005295105
$ws.Cells.Item(419, 4).Value = 46065.385636574072
$ws.Cells.Item(419, 5).Value = 'FRQ_ECO_SP_GUARULHOS_2'
$ws.Cells.Item(419, 6).Value = 'Detratores'
$ws.Cells.Item(419, 7).Value = 'não houve visita técnica.'
$ws.Cells.Item(419, 8).Value = 'Campo'
$ws.Cells.Item(419, 9).Value = 'Técnico não cumpriu a agenda'

$ws.Cells.Item(420, 1).Value = 419
$ws.Cells.Item(420, 2).Value = 'Pós OS'
$ws.Cells.Item(420, 3).Value = 8005303841
$ws.Cells.Item(420, 4).Value = 46065.392291666663
$ws.Cells.Item(420, 5).Value = 'FRQ_ECO_SP_OSASCO'
$ws.Cells.Item(420, 6).Value = 'Detratores'
$ws.Cells.Item(420, 7).Value = 'NÃO VIERAM'
$ws.Cells.Item(420, 8).Value = 'Campo'
$ws.Cells.Item(420, 9).Value = 'Técnico não cumpriu a agenda'

$ws.Cells.Item(421, 1).Value = 420
$ws.Cells.Item(421, 2).Value = 'Pós OS'
$ws.Cells.Item(421, 3).Value = 8005257619
$ws.Cells.Item(421, 4).Value = 46065.419849537036
$ws.Cells.Item(421, 5).Value = 'FRQ_ECO_RJ_OESTE'
$ws.Cells.Item(421, 6).Value = 'Detratores'
$ws.Cells.Item(421, 7).Value = 'Muito demorado agenda'
$ws.Cells.Item(421, 8).Value = 'Capacidade'
$ws.Cells.Item(421, 9).Value = 'Data Distante'

$ws.Cells.Item(422, 1).Value = 421
$ws.Cells.Item(422, 2).Value = 'Pós OS'
$ws.Cells.Item(422, 3).Value = 8005291791
$ws.Cells.Item(422, 4).Value = 46065.454837962963
$ws.Cells.Item(422, 5).Value = 'FRQ_ECO_SP_ZONA_SUL_03'
$ws.Cells.Item(422, 6).Value = 'Detratores'
$ws.Cells.Item(422, 7).Value = 'O purificador é bom quando funciona...Meu problema é que veio um técnico viu o problema, pediu a peça, mas demorou para vir a peça e quando veio precisa de mais uma...O filtro está parado faz 1 mês'
$ws.Cells.Item(422, 8).Value = 'Supply'
$ws.Cells.Item(422, 9).Value = 'Falta de estoque/peças'

$ws.Cells.Item(423, 1).Value = 422
$ws.Cells.Item(423, 2).Value = 'Pós OS'
$ws.Cells.Item(423, 3).Value = 8005237745
$ws.Cells.Item(423, 4).Value = 46065.505393518521
$ws.Cells.Item(423, 5).Value = 'FRQ_ECO_SP_S B CAMPO'
$ws.Cells.Item(423, 6).Value = 'Detratores'
$ws.Cells.Item(423, 7).Value = 'Não posso reclamar do purificador.. mas sim sobre assistência técnica..Vieram na minha casa constou que o purificador tinha que trocar .. foi agendado para depois de 1 mês (Janeiro).. ninguém veio e nem deu satisfação… liguei reagendaram para Fevereiro, fiquei praticamente 2 meses sem purificador .. mas o pagamento está em dia.. não descontaram os dias que tive que comprar água'
$ws.Cells.Item(423, 8).Value = 'Supply'
$ws.Cells.Item(423, 9).Value = 'Falta de estoque/peças'

$ws.Cells.Item(424, 1).Value = 423
$ws.Cells.Item(424, 2).Value = 'Pós OS'
$ws.Cells.Item(424, 3).Value = 8005285541
$ws.Cells.Item(424, 4).Value = 46065.511319444442
$ws.Cells.Item(424, 5).Value = 'FRQ_ECO_PE_RECIFE'
$ws.Cells.Item(424, 6).Value = 'Detratores'
$ws.Cells.Item(424, 7).Value = 'Muito ruim. Meu purificador ficou vazando água na bancada toda. To esperando nova vista agendada para amanhã'
$ws.Cells.Item(424, 8).Value = 'Qualidade do Produto'
$ws.Cells.Item(424, 9).Value = 'Vazamento'

$ws.Cells.Item(425, 1).Value = 424
$ws.Cells.Item(425, 2).Value = 'Pós OS'
$ws.Cells.Item(425, 3).Value = 8005303305
$ws.Cells.Item(425, 4).Value = 46065.514027777783
$ws.Cells.Item(425, 5).Value = 'FRQ_ECO_SP_CAMPINAS_2'
$ws.Cells.Item(425, 6).Value = 'Detratores'
$ws.Cells.Item(425, 7).Value = 'Acho um bom produto, mas caro. Pelo valor mensal , vc compra um filtro premium por ano.'
$ws.Cells.Item(425, 8).Value = 'Outros'
$ws.Cells.Item(425, 9).Value = 'Preço elevado'

$ws.Cells.Item(426, 1).Value = 425
$ws.Cells.Item(426, 2).Value = 'Pós OS'
$ws.Cells.Item(426, 3).Value = 8005241764
$ws.Cells.Item(426, 4).Value = 46065.536782407413
$ws.Cells.Item(426, 5).Value = 'FRQ_ECO_SP_OSASCO'
$ws.Cells.Item(426, 6).Value = 'Neutros'
$ws.Cells.Item(426, 7).Value = 'O purificador é ótimo. Mas o serviço de agendamento fica a desejar.'
$ws.Cells.Item(426, 8).Value = 'Outros'
$ws.Cells.Item(426, 9).Value = 'Satisfação geral'

$ws.Cells.Item(427, 1).Value = 426
$ws.Cells.Item(427, 2).Value = 'Pós OS'
$ws.Cells.Item(427, 3).Value = 8005283681
$ws.Cells.Item(427, 4).Value = 46065.55097222222
$ws.Cells.Item(427, 5).Value = 'AT_ECO_CE_FORTALEZA'
$ws.Cells.Item(427, 6).Value = 'Detratores'
$ws.Cells.Item(427, 7).Value = 'Meu purificador não está gelando.O técnico ficou de retornar e até o dia de hoje não retornou.Eu minha família estamos sem água gelada.'
$ws.Cells.Item(427, 8).Value = 'Campo'
$ws.Cells.Item(427, 9).Value = 'Reincidência'

$ws.Cells.Item(428, 1).Value = 427
$ws.Cells.Item(428, 2).Value = 'Pós OS'
$ws.Cells.Item(428, 3).Value = 8005285931
$ws.Cells.Item(428, 4).Value = 46065.721574074072
$ws.Cells.Item(428, 5).Value = 'FRQ_ECO_SP_ZONA_SUL_03'
$ws.Cells.Item(428, 6).Value = 'Neutros'
$ws.Cells.Item(428, 7).Value = 'A visita de revisão de 6 meses achei ruim, pois a técnica não veio com o filtro para substituição. Acreditávamos que o filtro seria substituído nesta revisão, mas ela informou que agora ia solicitar. Por que não vem com o filtro? Teremos que aguardar novo contato para novo agendamento.'
$ws.Cells.Item(428, 8).Value = 'Campo'
$ws.Cells.Item(428, 9).Value = 'Qualidade da manutenção'

$ws.Cells.Item(429, 1).Value = 428
$ws.Cells.Item(429, 2).Value = 'Pós OS'
$ws.Cells.Item(429, 3).Value = 8005259426
$ws.Cells.Item(429, 4).Value = 46065.837685185194
$ws.Cells.Item(429, 5).Value = 'FRQ_ECO_RJ_OESTE'
$ws.Cells.Item(429, 6).Value = 'Detratores'
$ws.Cells.Item(429, 7).Value = 'Não houve visita. Não cancelaram, não avisaram que não viriam. Achei total descaso.'
$ws.Cells.Item(429, 8).Value = 'Campo'
$ws.Cells.Item(429, 9).Value = 'Técnico não cumpriu a agenda'

$ws.Cells.Item(430, 1).Value = 429
$ws.Cells.Item(430, 2).Value = 'Pós OS'
$ws.Cells.Item(430, 3).Value = 8005279203
$ws.Cells.Item(430, 4).Value = 46065.852824074071
$ws.Cells.Item(430, 5).Value = 'FRQ_ECO_SP_OSASCO'
$ws.Cells.Item(430, 6).Value = 'Neutros'
$ws.Cells.Item(430, 7).Value = 'Eu gosto muito da água do purificador já uso  essa água a algum tempo, mais ultimamente está complicado estou tendo o mesmo problema a alguns meses já foi feita a troca do aparelho por 2 vezes sempre pelo mesmo motivo. E o problema continua a água só sai quente e agora está quente mesmo sai até fumaça e demora muito pra para um técnico vir ver o problema e depois mais tempão pra fazer a troca com isso já tenho uns 3 meses sem filtro e pelo jeito vou continuar pois o problema continua.Agua fervendo'
$ws.Cells.Item(430, 8).Value = 'Qualidade do Produto'
$ws.Cells.Item(430, 9).Value = 'Não gela'

$ws.Cells.Item(431, 1).Value = 430
$ws.Cells.Item(431, 2).Value = 'Instalação'
$ws.Cells.Item(431, 3).Value = 8005300765
$ws.Cells.Item(431, 4).Value = 46065.854409722233
$ws.Cells.Item(431, 5).Value = 'FRQ_ECO_SC_FLORIANOPOLIS2'
$ws.Cells.Item(431, 6).Value = 'Detratores'
$ws.Cells.Item(431, 7).Value = 'Boa noite!A pessoa responsável foi cordial e prestativa. Mas notamos que ficamos com um vazamento na torneira. Tentei contato pelo WhatsApp que fez a venda para programar o retorno da pessoa aqui e nao consegui mais contato. Ninguem me retorna. Inclusive gostaria de saber qual é o canal para solução de problemas. Obrigada'
$ws.Cells.Item(431, 8).Value = 'Campo'
$ws.Cells.Item(431, 9).Value = 'Qualidade da instalação'

$ws.Cells.Item(432, 1).Value = 431
$ws.Cells.Item(432, 2).Value = 'Pós OS'
$ws.Cells.Item(432, 3).Value = 8005285560
$ws.Cells.Item(432, 4).Value = 46066.37226851852
$ws.Cells.Item(432, 5).Value = 'FRQ_ECO_SP_ZONA_SUL_03'
$ws.Cells.Item(432, 6).Value = 'Detratores'
$ws.Cells.Item(432, 7).Value = 'O aparelho está velho, fazendo barulho e vcs dizem que é normal! A qualidade da prestação de serviço só cai!'
$ws.Cells.Item(432, 8).Value = 'Qualidade do Produto'
$ws.Cells.Item(432, 9).Value = 'Ruído'

$ws.Cells.Item(433, 1).Value = 432
$ws.Cells.Item(433, 2).Value = 'Pós OS'
$ws.Cells.Item(433, 3).Value = 8005280193
$ws.Cells.Item(433, 4).Value = 46066.40116898148
$ws.Cells.Item(433, 5).Value = 'FRQ_ECO_SP_CAMPINAS_2'
$ws.Cells.Item(433, 6).Value = 'Detratores'
$ws.Cells.Item(433, 7).Value = 'Não funciona a água com gás. O técnico esteve em casa,trocou umas peças, e nada de funcionar Péssimo atendimento'
$ws.Cells.Item(433, 8).Value = 'Qualidade do Produto'
$ws.Cells.Item(433, 9).Value = 'Funcionamento geral'

$ws.Cells.Item(434, 1).Value = 433
$ws.Cells.Item(434, 2).Value = 'Pós OS'
$ws.Cells.Item(434, 3).Value = 8005280801
$ws.Cells.Item(434, 4).Value = 46066.414884259262
$ws.Cells.Item(434, 5).Value = 'FRQ_ECO_SP_ZONA_SUL_03'
$ws.Cells.Item(434, 6).Value = 'Detratores'
$ws.Cells.Item(434, 7).Value = 'Existem dois contratos com a nossa empresa referente a dois purificadores, e quando o técnico veio aqui, fez a manutenção de apenas um, e solicitou para agendarmos outra visita para o segundo bebedouro conforme normas da Brastemp. A logística de vocês não era essa, pois quando o técnico vinha, realizava o procedimento nos dois equipamentos, evitando o transtorno e o custo de ter que voltar aqui outra vez, daqui 30 dias (prazo disponível na central). Nunca foi dessa forma, mas com a mudança de gestão, isso mudou para pior.'
$ws.Cells.Item(434, 8).Value = 'Capacidade'
$ws.Cells.Item(434, 9).Value = 'Agenda distante'

$ws.Cells.Item(435, 1).Value = 434
$ws.Cells.Item(435, 2).Value = 'Pós OS'
$ws.Cells.Item(435, 3).Value = 8005257620
$ws.Cells.Item(435, 4).Value = 46066.467534722222
$ws.Cells.Item(435, 5).Value = 'FRQ_ECO_SP_S B CAMPO'
$ws.Cells.Item(435, 6).Value = 'Detratores'
$ws.Cells.Item(435, 7).Value = 'Na verdade, minha insatisfação está com o prestador de serviços e/ou a empresa. Ele veio, e não tive nenhum problema de comportamento ou doisa semelhante. No entanto, ele identificou que o aparelho está com uma peça defeituosa. Informou que avisaria a empresa e que entrariam em contato comigo para agendar a substituição. Avisou que ficaria pingando até está troca e realmente ficou. No entanto, até o presente momento, ninguém entrou em contato comigo para falar sobre a substituição desta peça. Pago em dia o produto, mas a manutenção não está de acordo.'
$ws.Cells.Item(435, 8).Value = 'Campo'
$ws.Cells.Item(435, 9).Value = 'Reincidência'

$ws.Cells.Item(436, 1).Value = 435
$ws.Cells.Item(436, 2).Value = 'Pós OS'
$ws.Cells.Item(436, 3).Value = 8005282612
$ws.Cells.Item(436, 4).Value = 46066.502013888887
$ws.Cells.Item(436, 5).Value = 'FRQ_ECO_SP_ZONA_SUL_03'
$ws.Cells.Item(436, 6).Value = 'Detratores'
$ws.Cells.Item(436, 7).Value = 'Demorou muito para  vocês virem, fiquei 15 dias comprando água .'
$ws.Cells.Item(436, 8).Value = 'Capacidade'
$ws.Cells.Item(436, 9).Value = 'Agenda distante'

$ws.Cells.Item(437, 1).Value = 436
$ws.Cells.Item(437, 2).Value = 'Pós OS'
$ws.Cells.Item(437, 3).Value = 8005300267
$ws.Cells.Item(437, 4).Value = 46066.503460648149
$ws.Cells.Item(437, 5).Value = 'FRQ_ECO_SP_SJCAMPOS_3'
$ws.Cells.Item(437, 6).Value = 'Detratores'
$ws.Cells.Item(437, 7).Value = 'A pessoa que veio é ótima educada e prestativa, mas não resolveu o meu problema, continuou sem água fria'
$ws.Cells.Item(437, 8).Value = 'Qualidade do Produto'
$ws.Cells.Item(437, 9).Value = 'Não gela'

$ws.Cells.Item(438, 1).Value = 437
$ws.Cells.Item(438, 2).Value = 'Pós OS'
$ws.Cells.Item(438, 3).Value = 8005277871
$ws.Cells.Item(438, 4).Value = 46066.503738425927
$ws.Cells.Item(438, 5).Value = 'FRQ_ECO_BA_SALVADOR'
$ws.Cells.Item(438, 6).Value = 'Neutros'
$ws.Cells.Item(438, 7).Value = 'A comunicação entre o cliente/usuário e a brastemp é ruim, a cada atualização o aparelho fica ruim, perdendo funções importantes p a rotina de uma casa: perdeu a primeira função, a de ter a possibilidade de ter água em caso de falta de energia, depois perdeu a função, que era excelente, de timer (ele programava p desligar sozinho com 200ml, 500ml e 1l) essa função era excelente e providencial. Uso agora o modelo mais novo e os botões de acionamento são os piores de todas as versões que já tive (tive todas). Volta com o botão de timer, por favor.'
$ws.Cells.Item(438, 8).Value = 'Qualidade do Produto'
$ws.Cells.Item(438, 9).Value = 'Funcionamento geral'

$ws.Cells.Item(439, 1).Value = 438
$ws.Cells.Item(439, 2).Value = 'Pós OS'
$ws.Cells.Item(439, 3).Value = 8005295803
$ws.Cells.Item(439, 4).Value = 46066.503981481481
$ws.Cells.Item(439, 5).Value = 'FRQ_ECO_SC_JOINVILLE2'
$ws.Cells.Item(439, 6).Value = 'Detratores'
$ws.Cells.Item(439, 7).Value = 'Aparelho da problema desde o dia da instalação'
$ws.Cells.Item(439, 8).Value = 'Qualidade do Produto'
$ws.Cells.Item(439, 9).Value = 'Funcionamento geral'

$ws.Cells.Item(440, 1).Value = 439
$ws.Cells.Item(440, 2).Value = 'Pós OS'
$ws.Cells.Item(440, 3).Value = 8005280807
$ws.Cells.Item(440, 4).Value = 46066.555960648147
$ws.Cells.Item(440, 5).Value = 'FRQ_ECO_SP_ZONA_SUL_03'
$ws.Cells.Item(440, 6).Value = 'Detratores'
$ws.Cells.Item(440, 7).Value = 'O purificador estava cheio de mofo. Acredito que as manutenções anteriores tenham sido superficiais. Fato é que inclusive nesta, se não insistíssemos com o técnico ele não teria feita a manutenção correta. Lembrando que eu abri o chamado há quase 20 dias pois o cheio e gosto de ovo na água estavam insuportáveis. Impossibilitando o consumo. Aproveito para pedir o reembolso do dias que não pude utilizar o purificador e fiquei tendo que comprar águas no supermercado Meu custo foi superior a 100 reais.'
$ws.Cells.Item(440, 8).Value = 'Campo'
$ws.Cells.Item(440, 9).Value = 'Qualidade da manutenção'

$ws.Cells.Item(441, 1).Value = 440
$ws.Cells.Item(441, 2).Value = 'Pós OS'
$ws.Cells.Item(441, 3).Value = 8005303612
$ws.Cells.Item(441, 4).Value = 46066.560902777783
$ws.Cells.Item(441, 5).Value = 'FRQ_ECO_PE_RECIFE'
$ws.Cells.Item(441, 6).Value = 'Neutros'
$ws.Cells.Item(441, 7).Value = 'O técnico veio quarta-feira dia 11/02/26 realizou o serviço do vazamento, porém, hoje voltou a vazar de novo.'
$ws.Cells.Item(441, 8).Value = 'campo'
$ws.Cells.Item(441, 9).Value = 'Reincidência'

$excel.Goto($ws.Range("A424"))
$ws.Range("D428").Select()